$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column F, matching style of existing header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Timestamp values for rows 2-7
$ws.Range("F2").Value = "2021-10-05 10:52:46.307460"
$ws.Range("F3").Value = "2021-10-05 10:52:46.307473"
$ws.Range("F4").Value = "2021-10-05 10:52:46.307477"
$ws.Range("F5").Value = "2021-10-05 10:52:46.307480"
$ws.Range("F6").Value = "2021-10-05 10:52:46.307484"
$ws.Range("F7").Value = "2021-10-05 10:52:46.307487"
